$d = $word.ActiveDocument

$d.Content.Find.Execute("30-1=29", $true, $false, $false, $false, $false, $true, 1, $false, "93-87=6", 2)
$d.Content.Find.Execute("34-31=3", $true, $false, $false, $false, $false, $true, 1, $false, "28+64=92", 2)
$d.Content.Find.Execute("79+10=89", $true, $false, $false, $false, $false, $true, 1, $false, "78-54=24", 2)
$d.Content.Find.Execute("7+60=67", $true, $false, $false, $false, $false, $true, 1, $false, "18+62=80", 2)
$d.Content.Find.Execute("64-22=42", $true, $false, $false, $false, $false, $true, 1, $false, "16+55=71", 2)
$d.Content.Find.Execute("52-49=3", $true, $false, $false, $false, $false, $true, 1, $false, "83-20=63", 2)
$d.Content.Find.Execute("69-20=49", $true, $false, $false, $false, $false, $true, 1, $false, "41-17=24", 2)
$d.Content.Find.Execute("59-20=39", $true, $false, $false, $false, $false, $true, 1, $false, "20+5=25", 2)
$d.Content.Find.Execute("86-50=36", $true, $false, $false, $false, $false, $true, 1, $false, "63-57=6", 2)
$d.Content.Find.Execute("65-17=48", $true, $false, $false, $false, $false, $true, 1, $false, "8+71=79", 2)
$d.Content.Find.Execute("72+22=94", $true, $false, $false, $false, $false, $true, 1, $false, "95-42=53", 2)
$d.Content.Find.Execute("90+6=96", $true, $false, $false, $false, $false, $true, 1, $false, "23-6=17", 2)
$d.Content.Find.Execute("70+10=80", $true, $false, $false, $false, $false, $true, 1, $false, "20+24=44", 2)
$d.Content.Find.Execute("64+24=88", $true, $false, $false, $false, $false, $true, 1, $false, "99-68=31", 2)
$d.Content.Find.Execute("41-18=23", $true, $false, $false, $false, $false, $true, 1, $false, "92-33=59", 2)
$d.Content.Find.Execute("97-20=77", $true, $false, $false, $false, $false, $true, 1, $false, "78-12=66", 2)
$d.Content.Find.Execute("44+50=94", $true, $false, $false, $false, $false, $true, 1, $false, "75-66=9", 2)
$d.Content.Find.Execute("84-54=30", $true, $false, $false, $false, $false, $true, 1, $false, "1+58=59", 2)
$d.Content.Find.Execute("26+67=93", $true, $false, $false, $false, $false, $true, 1, $false, "22+7=29", 2)
$d.Content.Find.Execute("34+52=86", $true, $false, $false, $false, $false, $true, 1, $false, "0+55=55", 2)
$d.Content.Find.Execute("1+34=35", $true, $false, $false, $false, $false, $true, 1, $false, "43-5=38", 2)
$d.Content.Find.Execute("12+0=12", $true, $false, $false, $false, $false, $true, 1, $false, "90-43=47", 2)
$d.Content.Find.Execute("76-27=49", $true, $false, $false, $false, $false, $true, 1, $false, "4-4=0", 2)
$d.Content.Find.Execute("68-65=3", $true, $false, $false, $false, $false, $true, 1, $false, "25+15=40", 2)
$d.Content.Find.Execute("71-46=25", $true, $false, $false, $false, $false, $true, 1, $false, "54+2=56", 2)
$d.Content.Find.Execute("4+70=74", $true, $false, $false, $false, $false, $true, 1, $false, "81+1=82", 2)
$d.Content.Find.Execute("77-6=71", $true, $false, $false, $false, $false, $true, 1, $false, "11+48=59", 2)
$d.Content.Find.Execute("52-15=37", $true, $false, $false, $false, $false, $true, 1, $false, "99-24=75", 2)
$d.Content.Find.Execute("33-11=22", $true, $false, $false, $false, $false, $true, 1, $false, "28+48=76", 2)
$d.Content.Find.Execute("99-95=4", $true, $false, $false, $false, $false, $true, 1, $false, "51+47=98", 2)
$d.Content.Find.Execute("40+8=48", $true, $false, $false, $false, $false, $true, 1, $false, "97-24=73", 2)
$d.Content.Find.Execute("25+45=70", $true, $false, $false, $false, $false, $true, 1, $false, "23+21=44", 2)
$d.Content.Find.Execute("17+23=40", $true, $false, $false, $false, $false, $true, 1, $false, "43-32=11", 2)
$d.Content.Find.Execute("76-15=61", $true, $false, $false, $false, $false, $true, 1, $false, "0+14=14", 2)
$d.Content.Find.Execute("75-35=40", $true, $false, $false, $false, $false, $true, 1, $false, "71+4=75", 2)
$d.Content.Find.Execute("95-90=5", $true, $false, $false, $false, $false, $true, 1, $false, "58-28=30", 2)
$d.Content.Find.Execute("53+44=97", $true, $false, $false, $false, $false, $true, 1, $false, "37+1=38", 2)
$d.Content.Find.Execute("28+46=74", $true, $false, $false, $false, $false, $true, 1, $false, "58-1=57", 2)
$d.Content.Find.Execute("27-14=13", $true, $false, $false, $false, $false, $true, 1, $false, "19+26=45", 2)
$d.Content.Find.Execute("84-57=27", $true, $false, $false, $false, $false, $true, 1, $false, "58-43=15", 2)
$d.Content.Find.Execute("39-18=21", $true, $false, $false, $false, $false, $true, 1, $false, "47+12=59", 2)
$d.Content.Find.Execute("49+32=81", $true, $false, $false, $false, $false, $true, 1, $false, "46+27=73", 2)
$d.Content.Find.Execute("56-54=2", $true, $false, $false, $false, $false, $true, 1, $false, "71-15=56", 2)
$d.Content.Find.Execute("41-3=38", $true, $false, $false, $false, $false, $true, 1, $false, "25+38=63", 2)
$d.Content.Find.Execute("45+49=94", $true, $false, $false, $false, $false, $true, 1, $false, "72-26=46", 2)
$d.Content.Find.Execute("5+37=42", $true, $false, $false, $false, $false, $true, 1, $false, "33+58=91", 2)
$d.Content.Find.Execute("63-36=27", $true, $false, $false, $false, $false, $true, 1, $false, "21+64=85", 2)
$d.Content.Find.Execute("44-40=4", $true, $false, $false, $false, $false, $true, 1, $false, "1+92=93", 2)
$d.Content.Find.Execute("77-12=65", $true, $false, $false, $false, $false, $true, 1, $false, "5+75=80", 2)
$d.Content.Find.Execute("86-5=81", $true, $false, $false, $false, $false, $true, 1, $false, "55+28=83", 2)
$d.Content.Find.Execute("82-68=14", $true, $false, $false, $false, $false, $true, 1, $false, "83-74=9", 2)
$d.Content.Find.Execute("7-7=0", $true, $false, $false, $false, $false, $true, 1, $false, "88-24=64", 2)
$d.Content.Find.Execute("88-46=42", $true, $false, $false, $false, $false, $true, 1, $false, "31+62=93", 2)
$d.Content.Find.Execute("1+88=89", $true, $false, $false, $false, $false, $true, 1, $false, "25-24=1", 2)
$d.Content.Find.Execute("75-71=4", $true, $false, $false, $false, $false, $true, 1, $false, "83-75=8", 2)
$d.Content.Find.Execute("5+85=90", $true, $false, $false, $false, $false, $true, 1, $false, "46-26=20", 2)
$d.Content.Find.Execute("42-21=21", $true, $false, $false, $false, $false, $true, 1, $false, "16+14=30", 2)
$d.Content.Find.Execute("41+12=53", $true, $false, $false, $false, $false, $true, 1, $false, "82-71=11", 2)
$d.Content.Find.Execute("58+36=94", $true, $false, $false, $false, $false, $true, 1, $false, "16+9=25", 2)
$d.Content.Find.Execute("77-9=68", $true, $false, $false, $false, $false, $true, 1, $false, "41-16=25", 2)
$d.Content.Find.Execute("36-30=6", $true, $false, $false, $false, $false, $true, 1, $false, "10+71=81", 2)
$d.Content.Find.Execute("17+38=55", $true, $false, $false, $false, $false, $true, 1, $false, "16+0=16", 2)
$d.Content.Find.Execute("82-33=49", $true, $false, $false, $false, $false, $true, 1, $false, "13-8=5", 2)
$d.Content.Find.Execute("76-75=1", $true, $false, $false, $false, $false, $true, 1, $false, "8+80=88", 2)
$d.Content.Find.Execute("67+22=89", $true, $false, $false, $false, $false, $true, 1, $false, "95-39=56", 2)
$d.Content.Find.Execute("60+7=67", $true, $false, $false, $false, $false, $true, 1, $false, "52-28=24", 2)
$d.Content.Find.Execute("15+39=54", $true, $false, $false, $false, $false, $true, 1, $false, "50+25=75", 2)
$d.Content.Find.Execute("64-15=49", $true, $false, $false, $false, $false, $true, 1, $false, "90-6=84", 2)
$d.Content.Find.Execute("13+85=98", $true, $false, $false, $false, $false, $true, 1, $false, "29+44=73", 2)
$d.Content.Find.Execute("1+57=58", $true, $false, $false, $false, $false, $true, 1, $false, "15+50=65", 2)
$d.Content.Find.Execute("49+40=89", $true, $false, $false, $false, $false, $true, 1, $false, "20+41=61", 2)
$d.Content.Find.Execute("19+34=53", $true, $false, $false, $false, $false, $true, 1, $false, "4+52=56", 2)
$d.Content.Find.Execute("92-26=66", $true, $false, $false, $false, $false, $true, 1, $false, "99-33=66", 2)
$d.Content.Find.Execute("96+3=99", $true, $false, $false, $false, $false, $true, 1, $false, "41-30=11", 2)
$d.Content.Find.Execute("48+35=83", $true, $false, $false, $false, $false, $true, 1, $false, "39+14=53", 2)
$d.Content.Find.Execute("58+23=81", $true, $false, $false, $false, $false, $true, 1, $false, "20+66=86", 2)
$d.Content.Find.Execute("5+79=84", $true, $false, $false, $false, $false, $true, 1, $false, "93-47=46", 2)
$d.Content.Find.Execute("72-37=35", $true, $false, $false, $false, $false, $true, 1, $false, "22+7=29", 2)
$d.Content.Find.Execute("41+40=81", $true, $false, $false, $false, $false, $true, 1, $false, "24+74=98", 2)
$d.Content.Find.Execute("2+88=90", $true, $false, $false, $false, $false, $true, 1, $false, "76-26=50", 2)
$d.Content.Find.Execute("76-62=14", $true, $false, $false, $false, $false, $true, 1, $false, "29+38=67", 2)
$d.Content.Find.Execute("5+46=51", $true, $false, $false, $false, $false, $true, 1, $false, "38+57=95", 2)
$d.Content.Find.Execute("89-4=85", $true, $false, $false, $false, $false, $true, 1, $false, "6+0=6", 2)
$d.Content.Find.Execute("46+8=54", $true, $false, $false, $false, $false, $true, 1, $false, "29+49=78", 2)
$d.Content.Find.Execute("11+39=50", $true, $false, $false, $false, $false, $true, 1, $false, "5+76=81", 2)
$d.Content.Find.Execute("34-23=11", $true, $false, $false, $false, $false, $true, 1, $false, "78-73=5", 2)
$d.Content.Find.Execute("77-69=8", $true, $false, $false, $false, $false, $true, 1, $false, "91-5=86", 2)
$d.Content.Find.Execute("43+33=76", $true, $false, $false, $false, $false, $true, 1, $false, "81-70=11", 2)
$d.Content.Find.Execute("52+2=54", $true, $false, $false, $false, $false, $true, 1, $false, "72-32=40", 2)
$d.Content.Find.Execute("42-10=32", $true, $false, $false, $false, $false, $true, 1, $false, "39-20=19", 2)
$d.Content.Find.Execute("95-47=48", $true, $false, $false, $false, $false, $true, 1, $false, "10+15=25", 2)
$d.Content.Find.Execute("20+54=74", $true, $false, $false, $false, $false, $true, 1, $false, "16-8=8", 2)
$d.Content.Find.Execute("45+17=62", $true, $false, $false, $false, $false, $true, 1, $false, "23+34=57", 2)
$d.Content.Find.Execute("71-39=32", $true, $false, $false, $false, $false, $true, 1, $false, "25+20=45", 2)
$d.Content.Find.Execute("10+55=65", $true, $false, $false, $false, $false, $true, 1, $false, "93+3=96", 2)
$d.Content.Find.Execute("11+11=22", $true, $false, $false, $false, $false, $true, 1, $false, "6+88=94", 2)
$d.Content.Find.Execute("38+6=44", $true, $false, $false, $false, $false, $true, 1, $false, "1+50=51", 2)
$d.Content.Find.Execute("85+2=87", $true, $false, $false, $false, $false, $true, 1, $false, "3+9=12", 2)
$d.Content.Find.Execute("43+44=87", $true, $false, $false, $false, $false, $true, 1, $false, "99-74=25", 2)
$d.Content.Find.Execute("28-11=17", $true, $false, $false, $false, $false, $true, 1, $false, "90-75=15", 2)
